# Insert a new column before column C ("minuto") to hold the new "localia"
# field. This shifts the existing columns C:M (minuto..partido) one column
# to the right, becoming D:N, exactly as the target diff expects.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Insert()

# Header for the new column.
$ws.Range("C1").Value = "localia"

# Data rows: Villarreal (rows 2-22) played at home ("local"),
# Real Madrid (rows 23-34) played away ("visitante").
$ws.Range("C2:C22").Value = "local"
$ws.Range("C23:C34").Value = "visitante"
